# Update Efna2-Epha3.xlsx with new TPM-derived values.
#
# 1. Remove the three "ECs" sending-cluster rows (originally rows 2-4); the
#    remaining "FAPs"/"MuSCs" sending-cluster rows shift up to rows 2-7.
# 2. Recompute the numeric (TPM-derived) columns I/J/M/N/O/P/Q/R/S/T for the
#    surviving rows using the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: delete the old rows 2-4 (Sending cluster == "ECs")
$ws.Range("A2:A4").EntireRow.Delete() | Out-Null

# Step 2: write the updated TPM values into the now-shifted rows (2-7)
$data = @{
    2 = @{ I = 0.7554960962715589;  J = 0.7554960962715588;
           M = 0.003058333333333333; N = 0.009175000000000001;
           O = 0.0001379486413073712; P = 0.0001379486413073712;
           Q = 0.006163694658333334;  R = 0.055473251925;
           S = 0.0001042196599936845; T = 0.0001042196599936845 }
    3 = @{ I = 0.7554960962715589;  J = 0.7554960962715588;
           O = 0.9939610820947024;  P = 0.9939610820947024;
           S = 0.7509337173684021;  T = 0.750933717368402 }
    4 = @{ I = 0.7554960962715589;  J = 0.7554960962715588;
           O = 0.005900969263990248; P = 0.005900969263990248;
           S = 0.004458159243163087; T = 0.004458159243163086 }
    5 = @{ I = 0.2445039037284412;  J = 0.2445039037284411;
           M = 0.003058333333333333; N = 0.009175000000000001;
           O = 0.0001379486413073712; P = 0.0001379486413073712;
           Q = 0.001994778547222222;  R = 0.017953006925;
           S = 0.00003372898131368676; T = 0.00003372898131368676 }
    6 = @{ I = 0.2445039037284412;  J = 0.2445039037284411;
           O = 0.9939610820947024;  P = 0.9939610820947024;
           S = 0.2430273647263003;  T = 0.2430273647263003 }
    7 = @{ I = 0.2445039037284412;  J = 0.2445039037284411;
           O = 0.005900969263990248; P = 0.005900969263990248;
           Q = 0.08532977769166666; R = 0.767967999225;
           S = 0.001442810020827162; T = 0.001442810020827162 }
}

foreach ($r in $data.Keys) {
    foreach ($col in $data[$r].Keys) {
        $ws.Range("$col$r").Value = $data[$r][$col]
    }
}
